# Auto-generated edit script: applies the cell-value changes described by the diff
# across sheets ALC, ARM, BSM, CRP, CUL, LTW, WVR (GSM has no changes in this diff).

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 721.2
$ws.Range("I41").Value = 349
$ws.Range("K41").Value = 349
$ws.Range("M41").Value = 91
$ws.Range("H51").Value = 7737.3335
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 7737.3335
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 7737.3335
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -8705.333500000001
$ws.Range("H53").Value = 1371.4375
$ws.Range("I53").Value = 158.8
$ws.Range("J53").Value = 3392.5
$ws.Range("K53").Value = 158.8
$ws.Range("L53").Value = 3392.5
$ws.Range("M53").Value = 478.2
$ws.Range("N53").Value = -4666.5
$ws.Range("H86").Value = 4403.5293
$ws.Range("I86").Value = 3884
$ws.Range("J86").Value = 4767.2
$ws.Range("K86").Value = 3884
$ws.Range("L86").Value = 4767.2
$ws.Range("M86").Value = -2761
$ws.Range("N86").Value = -7013.2
$ws.Range("H89").Value = 4403.5293
$ws.Range("I89").Value = 3884
$ws.Range("J89").Value = 4767.2
$ws.Range("K89").Value = 19420
$ws.Range("L89").Value = 23836
$ws.Range("M89").Value = -13804
$ws.Range("N89").Value = -35068

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5364.5557
$ws.Range("I32").Value = 4350.125
$ws.Range("K32").Value = 4350.125
$ws.Range("M32").Value = -4063.125
$ws.Range("H74").Value = 5524.4
$ws.Range("I74").Value = 5524.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 5524.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -4650.4
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 5524.4
$ws.Range("I77").Value = 5524.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 27622
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -23254
$ws.Range("N77").ClearContents()
$ws.Range("H110").Value = 2495.7727
$ws.Range("I110").Value = 1694.7
$ws.Range("J110").Value = 10506.5
$ws.Range("K110").Value = 1694.7
$ws.Range("L110").Value = 10506.5
$ws.Range("M110").Value = 350.3
$ws.Range("N110").Value = -14596.5
$ws.Range("H122").Value = 4075.6667
$ws.Range("I122").Value = 3889.6365
$ws.Range("K122").Value = 11668.9095
$ws.Range("M122").Value = -9218.9095
$ws.Range("H132").Value = 3150.4
$ws.Range("I132").Value = 2017.2858
$ws.Range("K132").Value = 6051.857400000001
$ws.Range("M132").Value = -3521.857400000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3278.261
$ws.Range("I86").Value = 2766.2
$ws.Range("J86").Value = 4238.375
$ws.Range("K86").Value = 2766.2
$ws.Range("L86").Value = 4238.375
$ws.Range("M86").Value = -1643.2
$ws.Range("N86").Value = -6484.375
$ws.Range("H89").Value = 3278.261
$ws.Range("I89").Value = 2766.2
$ws.Range("J89").Value = 4238.375
$ws.Range("K89").Value = 13831
$ws.Range("L89").Value = 21191.875
$ws.Range("M89").Value = -8215
$ws.Range("N89").Value = -32423.875
$ws.Range("H94").Value = 2043.65
$ws.Range("I94").Value = 2087.4443
$ws.Range("K94").Value = 2087.4443
$ws.Range("M94").Value = -1636.4443
$ws.Range("H134").Value = 2624.3845
$ws.Range("I134").Value = 2146.6667
$ws.Range("K134").Value = 6440.000100000001
$ws.Range("M134").Value = -3905.000100000001

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2896.5
$ws.Range("I58").Value = 1488.5625
$ws.Range("K58").Value = 1488.5625
$ws.Range("M58").Value = -1285.5625
$ws.Range("H86").Value = 4789.8184
$ws.Range("I86").Value = 6253.8
$ws.Range("J86").Value = 3569.8333
$ws.Range("K86").Value = 6253.8
$ws.Range("L86").Value = 3569.8333
$ws.Range("M86").Value = -5130.8
$ws.Range("N86").Value = -5815.8333
$ws.Range("H89").Value = 4789.8184
$ws.Range("I89").Value = 6253.8
$ws.Range("J89").Value = 3569.8333
$ws.Range("K89").Value = 31269
$ws.Range("L89").Value = 17849.1665
$ws.Range("M89").Value = -25653
$ws.Range("N89").Value = -29081.1665
$ws.Range("H99").Value = 2200.147
$ws.Range("I99").Value = 1999.5
$ws.Range("J99").Value = 2378.5
$ws.Range("K99").Value = 1999.5
$ws.Range("L99").Value = 2378.5
$ws.Range("M99").Value = -501.5
$ws.Range("N99").Value = -5374.5
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H126").Value = 2200.147
$ws.Range("I126").Value = 1999.5
$ws.Range("J126").Value = 2378.5
$ws.Range("K126").Value = 5998.5
$ws.Range("L126").Value = 7135.5
$ws.Range("M126").Value = -3528.5
$ws.Range("N126").Value = -12075.5
$ws.Range("H132").Value = 4549
$ws.Range("I132").Value = 3424.4736
$ws.Range("K132").Value = 10273.4208
$ws.Range("M132").Value = -7743.4208
$ws.Range("H134").Value = 4319.2
$ws.Range("I134").Value = 2478.8333
$ws.Range("K134").Value = 7436.499899999999
$ws.Range("M134").Value = -4901.499899999999
$ws.Range("H136").Value = 2896.5
$ws.Range("I136").Value = 1488.5625
$ws.Range("K136").Value = 4465.6875
$ws.Range("M136").Value = -1915.6875

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1335530.6
$ws.Range("I5").Value = 803.1875
$ws.Range("K5").Value = 2409.5625
$ws.Range("M5").Value = -2297.5625
$ws.Range("H33").Value = 29.125
$ws.Range("I33").Value = 23.75
$ws.Range("J33").Value = 34.5
$ws.Range("K33").Value = 142.5
$ws.Range("L33").Value = 207
$ws.Range("M33").Value = 140.5
$ws.Range("N33").Value = -773
$ws.Range("H93").Value = 13006
$ws.Range("I93").Value = 10024
$ws.Range("J93").Value = 14000
$ws.Range("K93").Value = 30072
$ws.Range("L93").Value = 42000
$ws.Range("M93").Value = -28200
$ws.Range("N93").Value = -45744
$ws.Range("H131").Value = 6656793
$ws.Range("J131").Value = 4987645.5
$ws.Range("L131").Value = 14962936.5
$ws.Range("N131").Value = -14973016.5
$ws.Range("H132").Value = 4283.212
$ws.Range("I132").Value = 3045.158
$ws.Range("J132").Value = 5963.4287
$ws.Range("K132").Value = 27406.422
$ws.Range("L132").Value = 53670.85830000001
$ws.Range("M132").Value = -24876.422
$ws.Range("N132").Value = -58730.85830000001
$ws.Range("H135").Value = 1335530.6
$ws.Range("I135").Value = 803.1875
$ws.Range("K135").Value = 7228.6875
$ws.Range("M135").Value = -4693.6875

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 99235.5
$ws.Range("J6").Value = 99235.5
$ws.Range("L6").Value = 99235.5
$ws.Range("N6").Value = -99459.5
$ws.Range("H82").Value = 3367.04
$ws.Range("J82").Value = 3953.0557
$ws.Range("L82").Value = 3953.0557
$ws.Range("N82").Value = -4675.0557
$ws.Range("H85").Value = 3367.04
$ws.Range("J85").Value = 3953.0557
$ws.Range("L85").Value = 3953.0557
$ws.Range("N85").Value = -6449.0557
$ws.Range("H136").Value = 4611.4414
$ws.Range("I136").Value = 3643.76
$ws.Range("K136").Value = 10931.28
$ws.Range("M136").Value = -8381.280000000001

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("M88").ClearContents()
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("M91").ClearContents()
$ws.Range("N91").ClearContents()
$ws.Range("H126").Value = 2922.926
$ws.Range("I126").Value = 2377.3684
$ws.Range("J126").Value = 4218.625
$ws.Range("K126").Value = 7132.1052
$ws.Range("L126").Value = 12655.875
$ws.Range("M126").Value = -4662.1052
$ws.Range("N126").Value = -17595.875
$ws.Range("H132").Value = 3019.55
$ws.Range("I132").Value = 2178.2104
$ws.Range("K132").Value = 6534.6312
$ws.Range("M132").Value = -4004.6312
